$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the shopping list data: rotate the three item rows
# (onion/4, cabbage/3, apple/2) -> (cabbage/3, apple/2, onion/4)
$ws.Range("A2").Value = "cabbage"
$ws.Range("B2").Value = 3
$ws.Range("A3").Value = "apple"
$ws.Range("B3").Value = 2
$ws.Range("A4").Value = "onion"
$ws.Range("B4").Value = 4
